$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.919.07'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.834.02'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.20'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6933'
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07682'
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3044'
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.34'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07820'
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '93.26'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '1.835.85'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.096'
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6810'
$ws.Range("E15").Value = '  -1.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.587'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008263'
$ws.Range("E17").Value = '  -2.89%  '
$ws.Range("D18").Value = '28.923.62'
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.89'
$ws.Range("E19").Value = '  -3.16%  '
$ws.Range("D20").Value = '2.074.81'
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.68'
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.451'
$ws.Range("E23").Value = '  -2.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9995'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1506'
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.22'
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.749'
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.16'
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.541'
$ws.Range("E29").Value = '  -2.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.211'
$ws.Range("E30").Value = '  -1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.178'
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.194'
$ws.Range("E32").Value = '  -1.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05101'
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7807'
$ws.Range("E34").Value = '  +2.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.851'
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.143'
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.696'
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").Value = '1.294.05'
$ws.Range("E38").Value = '  +4.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01858'
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.704'
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9526'
$ws.Range("E41").Value = '  +5.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.148'
$ws.Range("E42").Value = '  +5.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.99'
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.680'
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000123'
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5166'
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("D48").Value = '1.975.06'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.81'
$ws.Range("E49").Value = '  -6.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.752'
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.956'
$ws.Range("E51").Value = '  -0.96%  '
